$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell E8 text: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Set active cell selection to E8
$ws.Activate() | Out-Null
$ws.Range("E8").Select() | Out-Null
